$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
}

Replace-Text "Google colab and Jupiter" "Google collab and Jupiter"
Replace-Text "will be the one updated as the model progresses" "will be updated as the model progresses"
Replace-Text "This file serves on purpose, to store the function F_InDF " "This file serves one purpose, to store the function F_InDF "
Replace-Text "In this function interventions and sampling strategies, and contamination scenario" "In this function, interventions, sampling strategies, and contamination scenario"
Replace-Text "If you run these one by one this will take days to run." "If you run this one by one this will take days to run."
Replace-Text "Chunk creating the output that will be used to create plots." "Chunk creates the output that will be used to create plots."
Replace-Text "Run lines 61-65: this will set up the iteration number we are running. " "Run lines 61-65: this will set up our iteration number. "
Replace-Text "Skip all the way to line 108 " "Skip to line 108 "
Replace-Text "Growth or die off during post" "Growth or die-off during post"
